$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp15"
$ws.Range("C2").Value = "Acvr2a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02770266666666667
$ws.Range("H2").Value = 0.083108
$ws.Range("I2").Value = 0.05180130905700151
$ws.Range("J2").Value = 0.05180130905700151
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 0.3967542953826667
$ws.Range("R2").Value = 3.570788658444
$ws.Range("S2").Value = 0.01527915445115451
$ws.Range("T2").Value = 0.01527915445115451

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp15"
$ws.Range("C3").Value = "Acvr2a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02770266666666667
$ws.Range("H3").Value = 0.083108
$ws.Range("I3").Value = 0.05180130905700151
$ws.Range("J3").Value = 0.05180130905700151
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.084169
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 0.7503037057506667
$ws.Range("R3").Value = 6.752733351756
$ws.Range("S3").Value = 0.02889447282323956
$ws.Range("T3").Value = 0.02889447282323956

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp15"
$ws.Range("C4").Value = "Acvr2a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02770266666666667
$ws.Range("H4").Value = 0.083108
$ws.Range("I4").Value = 0.05180130905700151
$ws.Range("J4").Value = 0.05180130905700151
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 0.1980682583408889
$ws.Range("R4").Value = 1.782614325068
$ws.Range("S4").Value = 0.007627681782607432
$ws.Range("T4").Value = 0.00762768178260743

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp15"
$ws.Range("C5").Value = "Acvr2a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.4126986666666667
$ws.Range("H5").Value = 1.238096
$ws.Range("I5").Value = 0.7717066171516261
$ws.Range("J5").Value = 0.7717066171516261
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 5.910621192858667
$ws.Range("R5").Value = 53.195590735728
$ws.Range("S5").Value = 0.2276202051469966
$ws.Range("T5").Value = 0.2276202051469966

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp15"
$ws.Range("C6").Value = "Acvr2a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.4126986666666667
$ws.Range("H6").Value = 1.238096
$ws.Range("I6").Value = 0.7717066171516261
$ws.Range("J6").Value = 0.7717066171516261
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 27.084169
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("Q6").Value = 11.17760043407467
$ws.Range("R6").Value = 100.598403906672
$ws.Range("S6").Value = 0.4304535210155654
$ws.Range("T6").Value = 0.4304535210155654

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp15"
$ws.Range("C7").Value = "Acvr2a"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.4126986666666667
$ws.Range("H7").Value = 1.238096
$ws.Range("I7").Value = 0.7717066171516261
$ws.Range("J7").Value = 0.7717066171516261
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 2.950708937512889
$ws.Range("R7").Value = 26.556380437616
$ws.Range("S7").Value = 0.113632890989064
$ws.Range("T7").Value = 0.113632890989064

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Bmp15"
$ws.Range("C8").Value = "Acvr2a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.09438566666666666
$ws.Range("H8").Value = 0.283157
$ws.Range("I8").Value = 0.1764920737913724
$ws.Range("J8").Value = 0.1764920737913724
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 1.351780286105666
$ws.Range("R8").Value = 12.166022574951
$ws.Range("S8").Value = 0.05205755808015543
$ws.Range("T8").Value = 0.05205755808015543

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Bmp15"
$ws.Range("C9").Value = "Acvr2a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.09438566666666666
$ws.Range("H9").Value = 0.283157
$ws.Range("I9").Value = 0.1764920737913724
$ws.Range("J9").Value = 0.1764920737913724
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.084169
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 2.556357347177667
$ws.Range("R9").Value = 23.007216124599
$ws.Range("S9").Value = 0.098446265596694
$ws.Range("T9").Value = 0.098446265596694

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Bmp15"
$ws.Range("C10").Value = "Acvr2a"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.09438566666666666
$ws.Range("H10").Value = 0.283157
$ws.Range("I10").Value = 0.1764920737913724
$ws.Range("J10").Value = 0.1764920737913724
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 0.6748377271385555
$ws.Range("R10").Value = 6.073539544247
$ws.Range("S10").Value = 0.02598825011452294
$ws.Range("T10").Value = 0.02598825011452293
